$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.481.59'
$ws.Range('E2').Value = '  +3.33%  '
$ws.Range('D3').Value = '1.602.51'
$ws.Range('E3').Value = '  +3.05%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('E6').Value = '  +7.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.91'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +10.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.35'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('E10').Value = '  +2.69%  '
$ws.Range('E11').Value = '  +2.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0913'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.20%  '
$ws.Range('D13').Value = '1.830.86'
$ws.Range('E13').Value = '  +2.99%  '
$ws.Range('D14').Value = '1.582.23'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').Value = '29.544.33'
$ws.Range('E15').Value = '  +3.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.535'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.74'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.54'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.21'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.42%  '
$ws.Range('E21').Value = '  +3.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  +3.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.07'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.64'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.29'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.65%  '
$ws.Range('E28').Value = '  +5.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.37'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0471'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.69%  '
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.24'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.41%  '
$ws.Range('D34').Value = '1.425.11'
$ws.Range('E34').Value = '  +2.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.10'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.50%  '
$ws.Range('E36').Value = '  -1.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.51'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.30'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0165'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.534'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.95'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '53.47'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +22.17%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.793'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.17%  '
$ws.Range('E46').Value = '  +1.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.86%  '
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('D49').Value = '1.742.36'
$ws.Range('E49').Value = '  +3.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.48'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.835'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.69%  '
